# Add a new column AT with data for 2024-10-24.
#
# Column A (header "台番号") has style s=1 for every row, and is used as the
# bulk style donor for AT. Rows that need the "hit" (yellow, s=2) or
# "big hit" (light-blue, s=3) highlight style instead borrow the format of
# an existing cell elsewhere in the same row that already carries that
# style, via Copy/PasteSpecial(xlPasteFormats) so the existing style indices
# in the workbook are reused instead of new ones being created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Give column AT the same column width as the other data columns (B..AS),
#    which is what produces the new <col .../> entry for column 46.
$ws.Columns("AT").ColumnWidth = $ws.Columns("B").ColumnWidth

# 2) Bulk-copy the plain style (s=1) from column A down the whole AT range.
$ws.Range("A1:A53").Copy()
$ws.Range("AT1:AT53").PasteSpecial(-4122)

# 3) Fix up the rows that need the yellow "hit" style (s=2): for each such
#    row, copy the format from a cell in that row that already has it.
$style2Donors = @{5="D"; 11="D"; 30="H"; 40="J"}
foreach ($r in $style2Donors.Keys) {
    $donor = $style2Donors[$r]
    $ws.Range("$donor$r").Copy()
    $ws.Range("AT$r").PasteSpecial(-4122)
}

# 4) Fix up the rows that need the blue "big hit" style (s=3) the same way.
$style3Donors = @{3="I"; 9="B"; 12="F"; 15="H"; 24="D"; 27="C"; 28="B"; 32="H"; 33="B"; 43="F"; 47="B"}
foreach ($r in $style3Donors.Keys) {
    $donor = $style3Donors[$r]
    $ws.Range("$donor$r").Copy()
    $ws.Range("AT$r").PasteSpecial(-4122)
}

# 5) Write the 52 numeric data values into AT2:AT53.
$values = @{2=211.8; 3=129.2; 4=163.2; 5=124.2; 6=148.8; 7=153.3; 8=160.1; 9=138.9; 10=329.6; 11=124.9; 12=136.7; 13=144.4; 14=246.4; 15=129.7; 16=169.3; 17=192.3; 18=297.6; 19=155.1; 20=187.4; 21=164.3; 22=151.7; 23=250.9; 24=125.2; 25=140; 26=201.3; 27=129.6; 28=134.7; 29=296.2; 30=123; 31=193.9; 32=136; 33=126.3; 34=146.8; 35=194.8; 36=140.4; 37=195.4; 38=181.2; 39=156.1; 40=122.3; 41=160.7; 42=159.1; 43=137.9; 44=151.8; 45=147; 46=204.8; 47=133.8; 48=140.5; 49=148.3; 50=179.6; 51=150.6; 52=144.6; 53=352.7}
foreach ($r in $values.Keys) {
    $ws.Range("AT$r").Value = $values[$r]
}

# 6) Write the header date into AT1 as literal text (like B1:AS1), not an
#    auto-converted date serial: force Text format, assign, then restore
#    the normal (General) style/format that the rest of row 1 uses.
$ws.Range("AT1").NumberFormat = "@"
$ws.Range("AT1").Value = "2024/10/24"
$ws.Range("A1").Copy()
$ws.Range("AT1").PasteSpecial(-4122)
